# Update Pertanggal 2 Januari 2025 18:51 WIB
#
# - Rename the "LookUp" sheet to "DataLookUp"
# - Make "DataLookUp" the active tab (was "MAIN")
# - Update the remembered selection on each sheet

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("MAIN")
$wsLookup = $wb.Worksheets.Item("LookUp")

# Rename the lookup sheet.
$wsLookup.Name = "DataLookUp"

# Move the selection on MAIN before we leave it (it keeps its own
# remembered selection/scroll position even once it's no longer active).
[void]$wsMain.Range("F13").Select()

# Switch the active tab to DataLookUp, with its own remembered selection.
[void]$wsLookup.Range("H13").Select()
$wsLookup.Activate()
